# Thbs2-Itga4 LR-pairs sheet: expand the Sending/Target cluster grid from
# {FAPs,sCs} x {ECs,sCs} (4 rows, no same-cluster pairs) to the full
# {ECs,FAPs,sCs} x {ECs,FAPs,sCs} grid (9 rows, including same-cluster
# pairs), with recomputed NATMI metric columns (E:T), per Dr Hou's advice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: full 3x3 Sending-cluster x Target-cluster grid for Thbs2 -> Itga4 (ECs, FAPs, sCs)
# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs2"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8911683333333333
$ws.Range("H2").Value = 2.673505
$ws.Range("I2").Value = 0.02693425114262819
$ws.Range("J2").Value = 0.02693425114262819
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.906497
$ws.Range("N2").Value = 68.719491
$ws.Range("O2").Value = 0.9446038650914245
$ws.Range("P2").Value = 0.9446038650914245
$ws.Range("Q2").Value = 20.413544753995
$ws.Range("R2").Value = 183.721902785955
$ws.Range("S2").Value = 0.0254421977326697
$ws.Range("T2").Value = 0.0254421977326697

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs2"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8911683333333333
$ws.Range("H3").Value = 2.673505
$ws.Range("I3").Value = 0.02693425114262819
$ws.Range("J3").Value = 0.02693425114262819
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1329193333333333
$ws.Range("N3").Value = 0.3987579999999999
$ws.Range("O3").Value = 0.005481244732096839
$ws.Range("P3").Value = 0.005481244732096839
$ws.Range("Q3").Value = 0.1184535007544444
$ws.Range("R3").Value = 1.06608150679
$ws.Range("S3").Value = 0.000147633222188504
$ws.Range("T3").Value = 0.000147633222188504

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Thbs2"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8911683333333333
$ws.Range("H4").Value = 2.673505
$ws.Range("I4").Value = 0.02693425114262819
$ws.Range("J4").Value = 0.02693425114262819
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.210428333333333
$ws.Range("N4").Value = 3.631285
$ws.Range("O4").Value = 0.04991489017647865
$ws.Range("P4").Value = 0.04991489017647865
$ws.Range("Q4").Value = 1.078695400436111
$ws.Range("R4").Value = 9.708258603925
$ws.Range("S4").Value = 0.001344420187769981
$ws.Range("T4").Value = 0.00134442018776998

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Thbs2"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.359699
$ws.Range("H5").Value = 73.07909699999999
$ws.Range("I5").Value = 0.7362360466408275
$ws.Range("J5").Value = 0.7362360466408276
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.906497
$ws.Range("N5").Value = 68.719491
$ws.Range("O5").Value = 0.9446038650914245
$ws.Range("P5").Value = 0.9446038650914245
$ws.Range("Q5").Value = 557.9953720644029
$ws.Range("R5").Value = 5021.958348579627
$ws.Range("S5").Value = 0.695451415276556
$ws.Range("T5").Value = 0.6954514152765561

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Thbs2"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.359699
$ws.Range("H6").Value = 73.07909699999999
$ws.Range("I6").Value = 0.7362360466408275
$ws.Range("J6").Value = 0.7362360466408276
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1329193333333333
$ws.Range("N6").Value = 0.3987579999999999
$ws.Range("O6").Value = 0.005481244732096839
$ws.Range("P6").Value = 0.005481244732096839
$ws.Range("Q6").Value = 3.237874951280665
$ws.Range("R6").Value = 29.14087456152599
$ws.Range("S6").Value = 0.004035489952229839
$ws.Range("T6").Value = 0.00403548995222984

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Thbs2"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.359699
$ws.Range("H7").Value = 73.07909699999999
$ws.Range("I7").Value = 0.7362360466408275
$ws.Range("J7").Value = 0.7362360466408276
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.210428333333333
$ws.Range("N7").Value = 3.631285
$ws.Range("O7").Value = 0.04991489017647865
$ws.Range("P7").Value = 0.04991489017647865
$ws.Range("Q7").Value = 29.48566986107166
$ws.Range("R7").Value = 265.3710287496449
$ws.Range("S7").Value = 0.03674914141204173
$ws.Range("T7").Value = 0.03674914141204173

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Thbs2"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.835938333333334
$ws.Range("H8").Value = 23.507815
$ws.Range("I8").Value = 0.2368297022165442
$ws.Range("J8").Value = 0.2368297022165442
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.906497
$ws.Range("N8").Value = 68.719491
$ws.Range("O8").Value = 0.9446038650914245
$ws.Range("P8").Value = 0.9446038650914245
$ws.Range("Q8").Value = 179.493897924685
$ws.Range("R8").Value = 1615.445081322165
$ws.Range("S8").Value = 0.2237102520821988
$ws.Range("T8").Value = 0.2237102520821988

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Thbs2"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.835938333333334
$ws.Range("H9").Value = 23.507815
$ws.Range("I9").Value = 0.2368297022165442
$ws.Range("J9").Value = 0.2368297022165442
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1329193333333333
$ws.Range("N9").Value = 0.3987579999999999
$ws.Range("O9").Value = 0.005481244732096839
$ws.Range("P9").Value = 0.005481244732096839
$ws.Range("Q9").Value = 1.041547699307778
$ws.Range("R9").Value = 9.373929293769999
$ws.Range("S9").Value = 0.001298121557678496
$ws.Range("T9").Value = 0.001298121557678496

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Thbs2"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.835938333333334
$ws.Range("H10").Value = 23.507815
$ws.Range("I10").Value = 0.2368297022165442
$ws.Range("J10").Value = 0.2368297022165442
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.210428333333333
$ws.Range("N10").Value = 3.631285
$ws.Range("O10").Value = 0.04991489017647865
$ws.Range("P10").Value = 0.04991489017647865
$ws.Range("Q10").Value = 9.484841776919446
$ws.Range("R10").Value = 85.363575992275
$ws.Range("S10").Value = 0.01182132857666695
$ws.Range("T10").Value = 0.01182132857666695

